# Apply crypto price/volume updates as described in the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number (e.g. "0.9998", "1.001").
# Excel would silently convert these to numeric values on assignment, but the source
# workbook stores them as text (inline strings), so force a text number format first.
$textFormatCells = @(
    "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15",
    "D16", "D17", "D18", "D20", "D23", "D24", "D25", "D26", "D27", "D29",
    "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D40", "D41",
    "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# New cell values (address -> value)
$updates = @{
    "D2" = "27.964.53"
    "E2" = "  -0.23%  "
    "D3" = "1.912.46"
    "E3" = "  +0.08%  "
    "D4" = "0.9998"
    "E4" = "  -0.44%  "
    "D5" = "313.15"
    "E5" = "  -1.14%  "
    "D7" = "0.5016"
    "E7" = "  +4.08%  "
    "D8" = "0.3813"
    "E8" = "  +0.10%  "
    "D9" = "0.07317"
    "E9" = "  -0.57%  "
    "D10" = "0.9122"
    "E10" = "  -2.41%  "
    "D11" = "21.24"
    "E11" = "  +2.00%  "
    "D12" = "0.07676"
    "E12" = "  -1.75%  "
    "D13" = "1.873.66"
    "E13" = "  -1.94%  "
    "D14" = "5.490"
    "E14" = "  -0.32%  "
    "D15" = "92.85"
    "E15" = "  +0.90%  "
    "D16" = "1.001"
    "E16" = "  -0.42%  "
    "D17" = "0.000008748"
    "E17" = "  -1.47%  "
    "D18" = "0.9994"
    "E18" = "  -0.42%  "
    "D19" = "27.979.32"
    "E19" = "  -0.27%  "
    "D20" = "14.67"
    "E20" = "  -0.61%  "
    "E21" = "  +0.28%  "
    "D22" = "2.117.30"
    "E22" = "  -1.36%  "
    "D23" = "10.86"
    "E23" = "  -0.45%  "
    "D24" = "6.625"
    "E24" = "  -0.23%  "
    "D25" = "153.22"
    "E25" = "  -2.44%  "
    "D26" = "1.848"
    "E26" = "  -3.30%  "
    "D27" = "2.204"
    "E27" = "  +3.40%  "
    "E28" = "  -0.34%  "
    "D29" = "115.58"
    "E29" = "  -1.39%  "
    "D30" = "4.930"
    "E30" = "  -1.00%  "
    "D31" = "0.09022"
    "E31" = "  +0.74%  "
    "D32" = "3.210"
    "E32" = "  -2.15%  "
    "D33" = "4.861"
    "E33" = "  +4.26%  "
    "E34" = "  -1.33%  "
    "D35" = "0.7814"
    "E35" = "  +0.83%  "
    "D36" = "2.623"
    "E36" = "  +0.13%  "
    "D37" = "0.02083"
    "E37" = "  +1.63%  "
    "D38" = "3.077"
    "E38" = "  +2.62%  "
    "E39" = "  -1.18%  "
    "D40" = "0.5551"
    "E40" = "  +0.22%  "
    "D41" = "0.05280"
    "E41" = "  -0.33%  "
    "D42" = "6.869"
    "E42" = "  -2.31%  "
    "D43" = "113.59"
    "E43" = "  +4.69%  "
    "D44" = "8.537"
    "E44" = "  +0.17%  "
    "D45" = "0.1519"
    "E45" = "  -0.70%  "
    "B46" = "EnergySwap"
    "C46" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D46" = "10.64"
    "E46" = "  -0.56%  "
    "B47" = "Decentraland"
    "C47" = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
    "D47" = "0.4828"
    "E47" = "  -0.16%  "
    "D48" = "0.9976"
    "E48" = "  -0.64%  "
    "D49" = "1.640"
    "E49" = "  -0.76%  "
    "D50" = "67.60"
    "E50" = "  -0.88%  "
    "D51" = "0.06055"
    "E51" = "  -0.29%  "
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
